# "after Sao Paulo 2010"
# Re-merge runs that had been split into multiple <a:r> elements back into a
# single run each (formatting was identical across the split runs, so this
# just collapses them, using the first run's rPr / the paragraph's pPr).

$p = $ppt.ActivePresentation

# --- Slide 11 ("Benefits" / "Details") -------------------------------------
$s11 = $p.Slides.Item(11)
$tr11 = $s11.Shapes.Item(2).TextFrame.TextRange

# "O" + "bject " + "oriented programming" -> "Object oriented programming"
$tr11.Characters(9, 27).Text = "Object oriented programming"

# "Add any method you " + "need" -> "Add any method you need"
$tr11.Characters(37, 23).Text = "Add any method you need"

# "Methods " + "are invoked with an appropriate node type, enforced during the " + "compilation."
# -> "Methods are invoked with an appropriate node type, enforced during the compilation."
$tr11.Characters(61, 83).Text = "Methods are invoked with an appropriate node type, enforced during the compilation."

# --- Slide 14 ("Object Lifecycle") ------------------------------------------
$s14 = $p.Slides.Item(14)
$tr14 = $s14.Shapes.Item(3).TextFrame.TextRange

# " " + "provides support for interacting with objects life " + "cycle."
# -> " provides support for interacting with objects life cycle."
$tr14.Characters(18, 58).Text = " provides support for interacting with objects life cycle."

# --- Slide 7 ("JCR Way") -----------------------------------------------------
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Item(1).TextFrame.TextRange

# "JCR Way " + "– Not type safe" -> "JCR Way – Not type safe"
$tr7.Characters(1, 23).Text = "JCR Way $([char]0x2013) Not type safe"
